$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "ACTIVE_METADATA"
$ws.Range("C8").Value = "{ name: [AA, AAA], value: [BB, BBB], weight: [0.5, 0.6] }"

$ws.Range("C8").Select()
